$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Duplicate row 6 (the original "Extension.valueBoolean" row) down to row 7,
#    preserving styles/formatting exactly, then retarget the path cells that
#    point at the parent element.
$ws.Range("A6:AJ6").Copy($ws.Range("A7:AJ7"))
$ws.Range("A7").Value = "Extension.value[x]"
$ws.Range("AE7").Value = "Extension.value[x]"

# 2. Rewrite row 6 in place to describe the value[x] slicing element instead
#    of the old valueBoolean leaf.
$ws.Range("B6").Value = ""
$ws.Range("E6").Value = "0"
$ws.Range("J6").Value = "boolean`n"
$ws.Range("K6").Value = "Value of extension"
$ws.Range("L6").Value = "Value of extension - may be a resource or one of a constrained set of the data types (see Extensibility in the spec for list)."
$ws.Range("AA6").Value = "type:`$this}`n"
$ws.Range("AB6").Value = ""
$ws.Range("AD6").Value = "closed"
$ws.Range("AE6").Value = "Extension.value[x]"

# 3. Extend the autofilter / filter-database range to include the new row.
#    (This also resets which rows are hidden, so hidden state is restored
#    explicitly afterwards.)
$ws.AutoFilterMode = $false
$ws.Range("A1:AJ7").AutoFilter(7, "<> ")
$ws.Range("A1:AJ7").AutoFilter(27, @(""), 7)

$name = $wb.Names.Item(1)
$name.RefersTo = "=Elements!`$A`$1:`$AJ`$7"

# 4. Grow the conditional formatting range by one row (still excluding the
#    last, newly-hidden detail row) and keep the same rules/dxfs.
$fc1 = $ws.Range("A2:AI5").FormatConditions.Item(1)
$fc1.ModifyAppliesToRange($ws.Range("A2:AI6"))

# 5. Column A shrinks slightly now that the widest text has changed.
$ws.Columns.Item(1).ColumnWidth = 18.166666666666668

# 6. Re-assert the hidden detail rows (rows 2-7 stay hidden; only the header
#    row 1 is visible), since toggling AutoFilterMode above redisplays them.
$ws.Rows.Item(2).Hidden = $true
$ws.Rows.Item(3).Hidden = $true
$ws.Rows.Item(4).Hidden = $true
$ws.Rows.Item(5).Hidden = $true
$ws.Rows.Item(6).Hidden = $true
$ws.Rows.Item(7).Hidden = $true
